$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Mes y año" metadata rows (column M) to use "null" like the other
# non-dimension columns (B, I, J), instead of duplicating column K's
# refPeriod metadata.
$ws.Range("M3").Value = "null"
$ws.Range("M4").Value = "null"
$ws.Range("M5").Value = "null"

# Update the view's selection / scroll position.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L10").Select()
